$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N17").Value = -2164.2501
$ws.Range("H17").Value = 609.4167
$ws.Range("L17").Value = 1828.2501
$ws.Range("J17").Value = 609.4167
$ws.Range("M47").Value = -5095
$ws.Range("N47").ClearContents()
$ws.Range("L47").Value = 0
$ws.Range("K47").Value = 6067
$ws.Range("J47").Value = 0
$ws.Range("I47").Value = 6067
$ws.Range("H47").Value = 6067
$ws.Range("I96").Value = 276.66666
$ws.Range("K96").Value = 829.9999799999999
$ws.Range("J96").Value = 447.5
$ws.Range("N96").Value = -4088.5
$ws.Range("L96").Value = 1342.5
$ws.Range("M96").Value = 543.0000200000001
$ws.Range("H96").Value = 319.375
$ws.Range("L113").Value = 3183.3333
$ws.Range("M113").Value = 449.6667000000002
$ws.Range("N113").Value = -9691.3333
$ws.Range("H113").Value = 2867.5
$ws.Range("K113").Value = 2804.3333
$ws.Range("J113").Value = 3183.3333
$ws.Range("I113").Value = 2804.3333
$ws.Range("M116").Value = 1675.3334
$ws.Range("I116").Value = 1766.6666
$ws.Range("L116").Value = 2800
$ws.Range("N116").Value = -9684
$ws.Range("H116").Value = 2283.3333
$ws.Range("K116").Value = 1766.6666
$ws.Range("J116").Value = 2800
$ws.Range("N127").Value = -14528.3999
$ws.Range("K127").Value = 1550.7273
$ws.Range("I127").Value = 516.9091
$ws.Range("L127").Value = 4608.3999
$ws.Range("J127").Value = 1536.1333
$ws.Range("M127").Value = 3409.2727
$ws.Range("H127").Value = 1104.9231
$ws.Range("J129").Value = 1002.53125
$ws.Range("H129").Value = 933.9737
$ws.Range("N129").Value = -13007.59375
$ws.Range("L129").Value = 3007.59375
$ws.Range("M137").Value = -18753802.5
$ws.Range("K137").Value = 18756352.5
$ws.Range("N137").Value = -17019.75
$ws.Range("I137").Value = 6252117.5
$ws.Range("H137").Value = 5002488.5
$ws.Range("J137").Value = 3973.25
$ws.Range("L137").Value = 11919.75
$ws.Range("N138").Value = -20855012
$ws.Range("J138").Value = 6948244
$ws.Range("K138").Value = 4814.1
$ws.Range("L138").Value = 20844732
$ws.Range("I138").Value = 1604.7
$ws.Range("M138").Value = 325.8999999999996
$ws.Range("H138").Value = 4905114.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3227
$ws.Range("I2").Value = 3867.5
$ws.Range("J2").Value = 2800
$ws.Range("N2").Value = -3026
$ws.Range("M2").Value = -3754.5
$ws.Range("K2").Value = 3867.5
$ws.Range("L2").Value = 2800
$ws.Range("N45").Value = -2596.8572
$ws.Range("L45").Value = 1842.8572
$ws.Range("M45").Value = -981.3478
$ws.Range("K45").Value = 1358.3478
$ws.Range("I45").Value = 1358.3478
$ws.Range("J45").Value = 1842.8572
$ws.Range("H45").Value = 1471.4
$ws.Range("M61").Value = -52684652
$ws.Range("H61").Value = 33400892
$ws.Range("J61").Value = 92215.09
$ws.Range("K61").Value = 52684864
$ws.Range("L61").Value = 92215.09
$ws.Range("N61").Value = -92639.09
$ws.Range("I61").Value = 52684864
$ws.Range("L63").Value = 4750
$ws.Range("N63").Value = -6122
$ws.Range("H63").Value = 3215
$ws.Range("I63").Value = 2873.889
$ws.Range("K63").Value = 2873.889
$ws.Range("M63").Value = -2187.889
$ws.Range("J63").Value = 4750
$ws.Range("M66").Value = -10937.445
$ws.Range("H66").Value = 3215
$ws.Range("K66").Value = 14369.445
$ws.Range("J66").Value = 4750
$ws.Range("N66").Value = -30614
$ws.Range("L66").Value = 23750
$ws.Range("I66").Value = 2873.889
$ws.Range("M74").Value = -11952995
$ws.Range("I74").Value = 11953869
$ws.Range("K74").Value = 11953869
$ws.Range("H74").Value = 9335405
$ws.Range("I77").Value = 11953869
$ws.Range("M77").Value = -59764977
$ws.Range("H77").Value = 9335405
$ws.Range("K77").Value = 59769345
$ws.Range("N88").Value = -9985.200000000001
$ws.Range("M88").Value = -2294
$ws.Range("K88").Value = 2700
$ws.Range("J88").Value = 9173.200000000001
$ws.Range("L88").Value = 9173.200000000001
$ws.Range("H88").Value = 8094.3335
$ws.Range("I88").Value = 2700
$ws.Range("H91").Value = 8094.3335
$ws.Range("K91").Value = 2700
$ws.Range("M91").Value = -1296
$ws.Range("L91").Value = 9173.200000000001
$ws.Range("J91").Value = 9173.200000000001
$ws.Range("N91").Value = -11981.2
$ws.Range("I91").Value = 2700
$ws.Range("I102").Value = 15874356
$ws.Range("L102").Value = 2699.4
$ws.Range("K102").Value = 15874356
$ws.Range("J102").Value = 2699.4
$ws.Range("N102").Value = -5943.4
$ws.Range("H102").Value = 10205908
$ws.Range("M102").Value = -15872734
$ws.Range("M116").Value = -1573.5
$ws.Range("I116").Value = 3867.5
$ws.Range("L116").Value = 2800
$ws.Range("N116").Value = -7388
$ws.Range("H116").Value = 3227
$ws.Range("K116").Value = 3867.5
$ws.Range("J116").Value = 2800
$ws.Range("I132").Value = 35334.793
$ws.Range("N132").Value = -230460.59
$ws.Range("H132").Value = 48902.547
$ws.Range("L132").Value = 225400.59
$ws.Range("K132").Value = 106004.379
$ws.Range("J132").Value = 75133.53
$ws.Range("M132").Value = -103474.379
$ws.Range("H136").Value = 33400892
$ws.Range("L136").Value = 276645.27
$ws.Range("J136").Value = 92215.09
$ws.Range("I136").Value = 52684864
$ws.Range("M136").Value = -158052042
$ws.Range("N136").Value = -281745.27
$ws.Range("K136").Value = 158054592

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I3").Value = 3867.5
$ws.Range("M3").Value = -3753.5
$ws.Range("K3").Value = 3867.5
$ws.Range("N3").Value = -3028
$ws.Range("H3").Value = 3227
$ws.Range("J3").Value = 2800
$ws.Range("L3").Value = 2800
$ws.Range("N86").Value = -17337.875
$ws.Range("M86").Value = -38560.332
$ws.Range("H86").Value = 25631.072
$ws.Range("L86").Value = 15091.875
$ws.Range("K86").Value = 39683.332
$ws.Range("I86").Value = 39683.332
$ws.Range("J86").Value = 15091.875
$ws.Range("K89").Value = 198416.66
$ws.Range("N89").Value = -86691.375
$ws.Range("M89").Value = -192800.66
$ws.Range("H89").Value = 25631.072
$ws.Range("J89").Value = 15091.875
$ws.Range("I89").Value = 39683.332
$ws.Range("L89").Value = 75459.375
$ws.Range("J94").Value = 957.1429000000001
$ws.Range("M94").Value = -378.9
$ws.Range("L94").Value = 957.1429000000001
$ws.Range("N94").Value = -1859.1429
$ws.Range("I94").Value = 829.9
$ws.Range("H94").Value = 882.2941
$ws.Range("K94").Value = 829.9
$ws.Range("L99").Value = 1090.3636
$ws.Range("M99").Value = 189.1765
$ws.Range("J99").Value = 1090.3636
$ws.Range("I99").Value = 1308.8235
$ws.Range("H99").Value = 1223
$ws.Range("K99").Value = 1308.8235
$ws.Range("N99").Value = -4086.3636
$ws.Range("N105").Value = -5727.3333
$ws.Range("H105").Value = 50001532
$ws.Range("I105").Value = 71429800
$ws.Range("J105").Value = 2233.3333
$ws.Range("K105").Value = 71429800
$ws.Range("L105").Value = 2233.3333
$ws.Range("M105").Value = -71428053
$ws.Range("K134").Value = 7424.3634
$ws.Range("M134").Value = -4889.3634
$ws.Range("I134").Value = 2474.7878
$ws.Range("H134").Value = 3242.3865
$ws.Range("L134").Value = 16635.5448
$ws.Range("J134").Value = 5545.1816
$ws.Range("N134").Value = -21705.5448

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("L31").Value = 6349.9
$ws.Range("J31").Value = 6349.9
$ws.Range("M31").Value = -1247.6666
$ws.Range("H31").Value = 2956.5588
$ws.Range("I31").Value = 1542.6666
$ws.Range("N31").Value = -6939.9
$ws.Range("K31").Value = 1542.6666
$ws.Range("H34").Value = 2956.5588
$ws.Range("I34").Value = 1542.6666
$ws.Range("M34").Value = -1340.6666
$ws.Range("L34").Value = 6349.9
$ws.Range("N34").Value = -6753.9
$ws.Range("J34").Value = 6349.9
$ws.Range("K34").Value = 1542.6666
$ws.Range("L99").Value = 3976.1667
$ws.Range("M99").Value = -1040.3572
$ws.Range("J99").Value = 3976.1667
$ws.Range("I99").Value = 2538.3572
$ws.Range("H99").Value = 3201.9614
$ws.Range("K99").Value = 2538.3572
$ws.Range("N99").Value = -6972.1667
$ws.Range("H126").Value = 3201.9614
$ws.Range("L126").Value = 11928.5001
$ws.Range("N126").Value = -16868.5001
$ws.Range("I126").Value = 2538.3572
$ws.Range("K126").Value = 7615.071599999999
$ws.Range("J126").Value = 3976.1667
$ws.Range("M126").Value = -5145.071599999999
$ws.Range("K134").Value = 7176.999899999999
$ws.Range("M134").Value = -4641.999899999999
$ws.Range("I134").Value = 2392.3333
$ws.Range("H134").Value = 34363.855
$ws.Range("L134").Value = 426803.25
$ws.Range("J134").Value = 142267.75
$ws.Range("N134").Value = -431873.25

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M39").ClearContents()
$ws.Range("K39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 2238.889
$ws.Range("L39").Value = 6716.667
$ws.Range("N39").Value = -7304.667
$ws.Range("H39").Value = 2238.889
$ws.Range("N104").Value = -16992.0001
$ws.Range("H104").Value = 3916.6667
$ws.Range("J104").Value = 3916.6667
$ws.Range("L104").Value = 11750.0001
$ws.Range("L113").Value = 1981.70271
$ws.Range("N113").Value = -6321.70271
$ws.Range("H113").Value = 578.1711
$ws.Range("J113").Value = 660.56757
$ws.Range("J131").Value = 1196.0667
$ws.Range("H131").Value = 1104.1569
$ws.Range("N131").Value = -13668.2001
$ws.Range("L131").Value = 3588.2001
$ws.Range("K134").Value = 3711.6921
$ws.Range("M134").Value = 1358.3079
$ws.Range("I134").Value = 1237.2307
$ws.Range("H134").Value = 3885.182
$ws.Range("L134").Value = 23130
$ws.Range("J134").Value = 7710
$ws.Range("N134").Value = -33270

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("K70").Value = 76282.14
$ws.Range("J70").Value = 5750
$ws.Range("M70").Value = -76012.14
$ws.Range("I70").Value = 76282.14
$ws.Range("N70").Value = -6290
$ws.Range("L70").Value = 5750
$ws.Range("H70").Value = 46893.75
$ws.Range("H73").Value = 46893.75
$ws.Range("J73").Value = 5750
$ws.Range("K73").Value = 76282.14
$ws.Range("M73").Value = -75346.14
$ws.Range("N73").Value = -7622
$ws.Range("L73").Value = 5750
$ws.Range("I73").Value = 76282.14
$ws.Range("H126").Value = 3000
$ws.Range("L126").Value = 9000
$ws.Range("N126").Value = -13940
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("J126").Value = 3000
$ws.Range("M126").ClearContents()
$ws.Range("I132").Value = 92658.17999999999
$ws.Range("N132").Value = -390718.25
$ws.Range("H132").Value = 107771.69
$ws.Range("L132").Value = 385658.25
$ws.Range("K132").Value = 277974.54
$ws.Range("J132").Value = 128552.75
$ws.Range("M132").Value = -275444.54

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K93").Value = 1425.8334
$ws.Range("I93").Value = 1425.8334
$ws.Range("H93").Value = 1688.8235
$ws.Range("M93").Value = -177.8334
